$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Update row 2: Application URL, Username, Password.
# Value2 avoids the automatic re-format (AutoCorrect / hyperlink recognition)
# that plain .Value can trigger on cells that already carry the Hyperlink style.
$ws.Range("A2").Value2 = "http://172.16.2.61:1616/UI#"
$ws.Range("B2").Value2 = "Administrator"
$ws.Range("C2").Value2 = "Tetherfi@930"

# Add hyperlink on C2, matching the mailto:<password> pattern used by the
# other password cells in this sheet.
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Tetherfi@930")

# Adding the hyperlink resets the cell format; restore the original
# "Hyperlink" look (vertically centred) so the cell style is unchanged.
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C2").VerticalAlignment = -4108

# Move selection to E10, matching the final cursor position in the workbook.
$ws.Range("E10").Select()
